$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Quantity updates (column F) for existing rows ---
$ws.Range("F12").Value = 100
$ws.Range("E15").Value = "660mW, 0,1%"
$ws.Range("F15").Value = 16
$ws.Range("F17").Value = 4

# New Packaging value on row 19 ("traversant" already exists as a shared string)
$ws.Range("C19").Value = "traversant"

$ws.Range("F24").Value = 10
$ws.Range("F32").Value = 10
$ws.Range("F38").Value = 5
$ws.Range("F40").Value = 2
$ws.Range("F41").Value = 4

# --- New components: push-buttons for MICHA v2 ---
# Fill columns A/B/C first (row by row), then column D, matching the
# shared-string insertion order of the authored change.
$ws.Range("A42").Value = "Bouton poussoir"
$ws.Range("B42").Value = "rouge, off - (on)"
$ws.Range("C42").Value = "Sur boitier"

$ws.Range("A43").Value = "Bouton poussoir"
$ws.Range("B43").Value = "vert, off - (on)"
$ws.Range("C43").Value = "Sur boitier"

$ws.Range("A44").Value = "Bouton poussoir"
$ws.Range("B44").Value = "jaune, off - (on)"
$ws.Range("C44").Value = "Sur boitier"

$ws.Range("D42").Value = "LP10A1AR"
$ws.Range("D43").Value = "LP10A1AG"
$ws.Range("D44").Value = "LP10A1AY"

# Ref column (D) for the new rows uses a centered, unlocked-cell style
$refCells = $ws.Range("D42:D44")
$refCells.HorizontalAlignment = -4108
$refCells.VerticalAlignment = -4108
$refCells.Locked = $false

$ws.Range("F42").Value = 2
$ws.Range("F43").Value = 2
$ws.Range("F44").Value = 2

$ws.Range("F44").Select()
